$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("158:158").Insert()

$ws.Cells.Item(158, 1).Value = 4
$ws.Cells.Item(158, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(158, 3).Value = "Los Lagos"
$ws.Cells.Item(158, 4).Value = 45258
$ws.Cells.Item(158, 5).Value = 10
$ws.Cells.Item(158, 6).Value = 100112026
$ws.Cells.Item(158, 7).Value = "Haba"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 180
$ws.Cells.Item(158, 11).Value = 15000
$ws.Cells.Item(158, 12).Value = 15000
$ws.Cells.Item(158, 13).Value = 15000
$ws.Cells.Item(158, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(158, 15).Value = "Región del Maule"
$ws.Cells.Item(158, 16).Value = 600
$ws.Cells.Item(158, 17).Value = 25
$ws.Cells.Item(158, 18).Value = "Hortaliza"
